# Refresh the scraped cryptocurrency price quotes (column D) on Sheet1.
#
# These cells hold the price as TEXT (e.g. "267.92"), not a number, so a
# plain Range.Value assignment of a numeric-looking string would get
# auto-converted to a real number by Excel (and also re-stamp the cell
# with a new/different style). To keep the cells as plain text with their
# original (unstyled) formatting untouched, each cell is set to a formula
# that evaluates to the literal text, then immediately replaced with its
# computed value via Copy / PasteSpecial(values only).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = "267.73"
    "D3"  = "22.96"
    "D4"  = "6.313"
    "D5"  = "0.06206"
    "D7"  = "6.694"
    "D8"  = "1.391"
    "D9"  = "0.8320"
    "D10" = "0.01366"
    "D11" = "0.1604"
    "D12" = "0.08306"
    "D13" = "0.03405"
    "D15" = "0.09284"
    "D16" = "3.898"
    "D18" = "0.04839"
    "D19" = "0.006252"
    "D20" = "0.005380"
    "D21" = "0.001090"
    "D22" = "0.0001500"
    "D23" = "3.771"
    "D24" = "2.368"
    "D25" = "0.3348"
    "D40" = "0.04674"
    "D41" = "0.006911"
    "D42" = "0.1155"
    "D43" = "0.003459"
    "D45" = "0.00006220"
    "D47" = "0.6999"
    "D48" = "0.1723"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $updates[$addr] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

$excel.CutCopyMode = 0
